# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the zh-cn and de-de
# sheets to reflect the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 02:43:48"
$wsZhCn.Range("H2").Value = "2016-03-14 02:44:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 02:43:51"
$wsDeDe.Range("H2").Value = "2016-03-14 02:44:09"
